# Logged Week 16 and performed season sim from Week 17
# Update target depth data totals on both the OFF and DEF sheets (row 2 = "H")

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 288
$wsOff.Range("C2").Value = 189
$wsOff.Range("D2").Value = 76
$wsOff.Range("E2").Value = 32

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 399
$wsDef.Range("C2").Value = 291
$wsDef.Range("D2").Value = 85
$wsDef.Range("E2").Value = 45
$wsDef.Range("F2").Value = 4
$wsDef.Range("G2").Value = 4
